# ML model retrained with all data
# Update predicted-probability columns J and K for every data row (1-51).
# Final state: column J == 1, column K == 0.6 for all rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1:J51").Value = 1
$ws.Range("K1:K51").Value = 0.6

# Restore the scroll position / selection recorded in the sheet view:
# window scrolled down so row 32 is the top-visible row, and the
# previously-selected K1 cell is extended to the full K1:K51 range.
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 5
$ws.Range("K1:K51").Select() | Out-Null
